$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-4: add K (=2), L (=999), M (=999)
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 11).Value = 2
    $ws.Cells.Item($r, 12).Value = 999
    $ws.Cells.Item($r, 13).Value = 999
}

# Rows 5-25: K already set; add L (=999), M (=999)
foreach ($r in 5..25) {
    $ws.Cells.Item($r, 12).Value = 999
    $ws.Cells.Item($r, 13).Value = 999
}

# Sheet view changes: selection moves from B2 to K4
$ws.Activate()
$ws.Range("K4").Select()
